$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new columns: "start time" (F) and "end time" (G) ---
# Set the number format to Text ("@") BEFORE writing the date-looking
# strings so Excel stores them as literal text (matching the source data,
# e.g. "2017/01/01") instead of auto-converting them into date serials.
$ws.Range("F1:G4").NumberFormat = "@"

# Header row
$ws.Range("F1").Value = "start time"
$ws.Range("G1").Value = "end time"

# Row 2 - course C01 / course-v1:FCUx+QA76+19004
$ws.Range("F2").Value = "2017/01/01"
$ws.Range("G2").Value = "2020/02/01"

# Row 3 - course C01 / course-v1:FCUx+QA76+19004 (admin test row)
$ws.Range("F3").Value = "2017/01/01"
$ws.Range("G3").Value = "2020/02/01"

# Row 4 - course C02 / DYUx/dyu10401/201511
$ws.Range("F4").Value = "2017/01/01"
$ws.Range("G4").Value = "2020/08/01"

# Move the active selection cursor (cosmetic, matches author's saved view)
$ws.Range("K7").Select() | Out-Null
